$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 148.57143
$ws.Range("I9").Value = 170
$ws.Range("J9").Value = 95
$ws.Range("K9").Value = 170
$ws.Range("L9").Value = 95
$ws.Range("M9").Value = -1
$ws.Range("N9").Value = -433
$ws.Range("H51").Value = 4962.6665
$ws.Range("J51").Value = 4962.6665
$ws.Range("L51").Value = 4962.6665
$ws.Range("N51").Value = -5930.6665
$ws.Range("H69").Value = 3197.6
$ws.Range("J69").Value = 3243.75
$ws.Range("L69").Value = 9731.25
$ws.Range("N69").Value = -11479.25
$ws.Range("H72").Value = 3197.6
$ws.Range("J72").Value = 3243.75
$ws.Range("L72").Value = 29193.75
$ws.Range("N72").Value = -37929.75
$ws.Range("H98").Value = 7016.4
$ws.Range("I98").Value = 7544.75
$ws.Range("J98").Value = 4903
$ws.Range("K98").Value = 7544.75
$ws.Range("L98").Value = 4903
$ws.Range("M98").Value = -6046.75
$ws.Range("N98").Value = -7899
$ws.Range("H112").Value = 1879.4615
$ws.Range("I112").Value = 1100
$ws.Range("J112").Value = 1910.64
$ws.Range("K112").Value = 3300
$ws.Range("L112").Value = 5731.92
$ws.Range("M112").Value = -2192
$ws.Range("N112").Value = -7947.92
$ws.Range("H122").Value = 7016.4
$ws.Range("I122").Value = 7544.75
$ws.Range("J122").Value = 4903
$ws.Range("K122").Value = 22634.25
$ws.Range("L122").Value = 14709
$ws.Range("M122").Value = -20184.25
$ws.Range("N122").Value = -19609
$ws.Range("H137").Value = 1011
$ws.Range("I137").Value = 851.5172
$ws.Range("K137").Value = 2554.5516
$ws.Range("M137").Value = -4.55159999999978
$ws.Range("H138").Value = 1281.46
$ws.Range("I138").Value = 836.7742
$ws.Range("J138").Value = 1481.2463
$ws.Range("K138").Value = 2510.3226
$ws.Range("L138").Value = 4443.7389
$ws.Range("M138").Value = 2629.6774
$ws.Range("N138").Value = -14723.7389
$ws.Range("H141").Value = 1215.3
$ws.Range("I141").Value = 993
$ws.Range("J141").Value = 2104.5
$ws.Range("K141").Value = 2979
$ws.Range("L141").Value = 6313.5
$ws.Range("M141").Value = 2201
$ws.Range("N141").Value = -16673.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1119.95
$ws.Range("I45").Value = 1145.2778
$ws.Range("K45").Value = 1145.2778
$ws.Range("M45").Value = -768.2778000000001
$ws.Range("H61").Value = 1155.3572
$ws.Range("I61").Value = 982.2273
$ws.Range("J61").Value = 1790.1666
$ws.Range("K61").Value = 982.2273
$ws.Range("L61").Value = 1790.1666
$ws.Range("M61").Value = -770.2273
$ws.Range("N61").Value = -2214.1666
$ws.Range("H74").Value = 1166.1111
$ws.Range("I74").Value = 749.5
$ws.Range("J74").Value = 1999.3334
$ws.Range("K74").Value = 749.5
$ws.Range("L74").Value = 1999.3334
$ws.Range("M74").Value = 124.5
$ws.Range("N74").Value = -3747.3334
$ws.Range("H77").Value = 1166.1111
$ws.Range("I77").Value = 749.5
$ws.Range("J77").Value = 1999.3334
$ws.Range("K77").Value = 3747.5
$ws.Range("L77").Value = 9996.667
$ws.Range("M77").Value = 620.5
$ws.Range("N77").Value = -18732.667
$ws.Range("H132").Value = 2793.3157
$ws.Range("I132").Value = 2404.6924
$ws.Range("K132").Value = 7214.0772
$ws.Range("M132").Value = -4684.0772
$ws.Range("H136").Value = 1155.3572
$ws.Range("I136").Value = 982.2273
$ws.Range("J136").Value = 1790.1666
$ws.Range("K136").Value = 2946.6819
$ws.Range("L136").Value = 5370.4998
$ws.Range("M136").Value = -396.6819
$ws.Range("N136").Value = -10470.4998

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5294.3228
$ws.Range("I134").Value = 751.5714
$ws.Range("K134").Value = 2254.7142
$ws.Range("M134").Value = 280.2857999999997

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1186.7106
$ws.Range("I31").Value = 880.1667
$ws.Range("J31").Value = 2336.25
$ws.Range("K31").Value = 880.1667
$ws.Range("L31").Value = 2336.25
$ws.Range("M31").Value = -585.1667
$ws.Range("N31").Value = -2926.25
$ws.Range("H34").Value = 1186.7106
$ws.Range("I34").Value = 880.1667
$ws.Range("J34").Value = 2336.25
$ws.Range("K34").Value = 880.1667
$ws.Range("L34").Value = 2336.25
$ws.Range("M34").Value = -678.1667
$ws.Range("N34").Value = -2740.25
$ws.Range("H58").Value = 1130.1515
$ws.Range("I58").Value = 825.7407
$ws.Range("J58").Value = 2500
$ws.Range("K58").Value = 825.7407
$ws.Range("L58").Value = 2500
$ws.Range("M58").Value = -622.7407
$ws.Range("N58").Value = -2906
$ws.Range("H122").Value = 730.6667
$ws.Range("I122").Value = 633.86957
$ws.Range("K122").Value = 1901.60871
$ws.Range("M122").Value = 548.39129
$ws.Range("H132").Value = 4109.675
$ws.Range("I132").Value = 4363.2
$ws.Range("K132").Value = 13089.6
$ws.Range("M132").Value = -10559.6
$ws.Range("H134").Value = 1892.5333
$ws.Range("I134").Value = 1919.0435
$ws.Range("J134").Value = 1805.4286
$ws.Range("K134").Value = 5757.1305
$ws.Range("L134").Value = 5416.2858
$ws.Range("M134").Value = -3222.1305
$ws.Range("N134").Value = -10486.2858
$ws.Range("H136").Value = 1130.1515
$ws.Range("I136").Value = 825.7407
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 2477.2221
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = 72.77790000000005
$ws.Range("N136").Value = -12600

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 319.9
$ws.Range("I41").Value = 365.33334
$ws.Range("J41").Value = 300.42856
$ws.Range("K41").Value = 1096.00002
$ws.Range("L41").Value = 901.28568
$ws.Range("M41").Value = -758.00002
$ws.Range("N41").Value = -1577.28568

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 5000
$ws.Range("J19").Value = 5000
$ws.Range("L19").Value = 5000
$ws.Range("N19").Value = -5576
$ws.Range("H113").Value = 1691.6666
$ws.Range("I113").Value = 1171.4286
$ws.Range("K113").Value = 1171.4286
$ws.Range("M113").Value = 998.5714
$ws.Range("H132").Value = 2997.7856
$ws.Range("I132").Value = 2563.4443
$ws.Range("J132").Value = 3779.6
$ws.Range("K132").Value = 7690.3329
$ws.Range("L132").Value = 11338.8
$ws.Range("M132").Value = -5160.3329
$ws.Range("N132").Value = -16398.8
$ws.Range("H134").Value = 32326
$ws.Range("J134").Value = 32326
$ws.Range("L134").Value = 96978
$ws.Range("N134").Value = -102048
$ws.Range("H135").Value = 39999
$ws.Range("J135").Value = 39999
$ws.Range("L135").Value = 39999
$ws.Range("N135").Value = -50139

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 146.56667
$ws.Range("I55").Value = 70.42105
$ws.Range("J55").Value = 278.0909
$ws.Range("K55").Value = 70.42105
$ws.Range("L55").Value = 278.0909
$ws.Range("M55").Value = 102.57895
$ws.Range("N55").Value = -624.0908999999999
$ws.Range("H132").Value = 20116.834
$ws.Range("I132").Value = 1325.8148
$ws.Range("J132").Value = 38907.85
$ws.Range("K132").Value = 3977.4444
$ws.Range("L132").Value = 116723.55
$ws.Range("M132").Value = -1447.4444
$ws.Range("N132").Value = -121783.55
$ws.Range("H136").Value = 5144.448
$ws.Range("I136").Value = 6904.5557
$ws.Range("J136").Value = 2264.2727
$ws.Range("K136").Value = 20713.6671
$ws.Range("L136").Value = 6792.8181
$ws.Range("M136").Value = -18163.6671
$ws.Range("N136").Value = -11892.8181

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 300
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 300
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H124").Value = 60214.5
$ws.Range("J124").Value = 60214.5
$ws.Range("L124").Value = 60214.5
$ws.Range("N124").Value = -70034.5
$ws.Range("H132").Value = 2358.1365
$ws.Range("I132").Value = 2013.7407
$ws.Range("K132").Value = 6041.2221
$ws.Range("M132").Value = -3511.2221
$ws.Range("H136").Value = 656.0455
$ws.Range("I136").Value = 649.6667
$ws.Range("K136").Value = 1949.0001
$ws.Range("M136").Value = 600.9999
